# Apply updated PLS predicted-variables matrix values to A2:B67
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 66,2
$arr[0,0] = -0.2285456021302903
$arr[0,1] = -0.9492488821387489
$arr[1,0] = -0.4473314759756136
$arr[1,1] = -1.056489456487291
$arr[2,0] = -0.498877152018381
$arr[2,1] = -1.530416122729806
$arr[3,0] = -0.1816220354970723
$arr[3,1] = -0.8217912260128744
$arr[4,0] = -0.1072946438059379
$arr[4,1] = 0.6712711822651735
$arr[5,0] = -0.3141801510651895
$arr[5,1] = -0.9402540617027475
$arr[6,0] = -0.4543624128096752
$arr[6,1] = -0.7531652619036333
$arr[7,0] = -0.5765577983985608
$arr[7,1] = -1.502978382799211
$arr[8,0] = -0.4572847029734084
$arr[8,1] = -0.5494389262496754
$arr[9,0] = -0.3152550681314725
$arr[9,1] = -1.153662761156768
$arr[10,0] = -0.2527015666617348
$arr[10,1] = -0.788308836981051
$arr[11,0] = -0.0749596063598424
$arr[11,1] = -0.2016235787603602
$arr[12,0] = -0.2587200422907012
$arr[12,1] = -0.2975752239844328
$arr[13,0] = -0.1330218959055534
$arr[13,1] = 0.03377220449208486
$arr[14,0] = -0.2021974204587239
$arr[14,1] = -0.2809066531631849
$arr[15,0] = 0.05013896101029428
$arr[15,1] = 0.4997205414402155
$arr[16,0] = 0.02251671216509762
$arr[16,1] = 0.3596779905058352
$arr[17,0] = 0.02259866784861345
$arr[17,1] = 0.4354922473442499
$arr[18,0] = -0.1430035481804015
$arr[18,1] = -0.02307601288939132
$arr[19,0] = 0.001160352665126801
$arr[19,1] = 0.2765955462398949
$arr[20,0] = 0.04358583917243557
$arr[20,1] = 0.673526736293705
$arr[21,0] = 0.00924777846651384
$arr[21,1] = 0.3464234081472974
$arr[22,0] = 0.7418093705544981
$arr[22,1] = 1.634204086777708
$arr[23,0] = 0.1348365744491544
$arr[23,1] = 0.728072978469803
$arr[24,0] = 0.1488531563726015
$arr[24,1] = 0.6985023472566052
$arr[25,0] = 0.1025546604068538
$arr[25,1] = 0.7177211418854293
$arr[26,0] = 0.2978774132295715
$arr[26,1] = 0.9750705587986496
$arr[27,0] = 0.5900573471219441
$arr[27,1] = 1.697635037468139
$arr[28,0] = 0.2003068906414907
$arr[28,1] = 0.7729900924801829
$arr[29,0] = 0.03147667668032685
$arr[29,1] = 0.594698810225738
$arr[30,0] = 0.1674018562183917
$arr[30,1] = 0.9243527168147343
$arr[31,0] = 0.116198597946313
$arr[31,1] = 0.8043329125191326
$arr[32,0] = 0.08896254713267385
$arr[32,1] = 0.4298682362414687
$arr[33,0] = 0.3941985824917741
$arr[33,1] = 1.37037824345451
$arr[34,0] = 0.2675368268508296
$arr[34,1] = 1.033297314844164
$arr[35,0] = 0.05768397769491911
$arr[35,1] = 0.3144763238187946
$arr[36,0] = 0.3391397541901515
$arr[36,1] = 1.761767527176945
$arr[37,0] = -0.07508363684379435
$arr[37,1] = -0.1281075718919199
$arr[38,0] = 0.1429002900458812
$arr[38,1] = 0.9651662496716948
$arr[39,0] = -0.0881347708681114
$arr[39,1] = 0.6563022605677905
$arr[40,0] = 0.2431616990283887
$arr[40,1] = 1.327993780115742
$arr[41,0] = 0.1951268474735019
$arr[41,1] = 1.168210126878793
$arr[42,0] = -0.1390907336765805
$arr[42,1] = -0.1155862812699684
$arr[43,0] = -0.09659956304960367
$arr[43,1] = 0.2157738959363079
$arr[44,0] = -0.1909791062870775
$arr[44,1] = -0.4998771490212042
$arr[45,0] = -0.1889153857511643
$arr[45,1] = -0.4661730655458792
$arr[46,0] = -0.2258045536236594
$arr[46,1] = -0.5174113087010024
$arr[47,0] = -0.221023749255509
$arr[47,1] = -0.5682130542313183
$arr[48,0] = -0.1560859975108737
$arr[48,1] = -0.3589638811933314
$arr[49,0] = -0.2108553482898883
$arr[49,1] = -0.6256082688301617
$arr[50,0] = -0.2108553482898883
$arr[50,1] = -0.6256082688301617
$arr[51,0] = -0.2003845824100484
$arr[51,1] = -0.4252332136001163
$arr[52,0] = -0.2031094970425705
$arr[52,1] = -0.5237129288572019
$arr[53,0] = -0.1714564274003043
$arr[53,1] = -0.3146240846005753
$arr[54,0] = -0.1577884220445736
$arr[54,1] = -0.4062569261786673
$arr[55,0] = -0.204320210983713
$arr[55,1] = -0.5218910148814445
$arr[56,0] = -0.1777825254214953
$arr[56,1] = -0.5904498575473008
$arr[57,0] = -0.2146862012650751
$arr[57,1] = -0.6542169155027874
$arr[58,0] = -0.2057395778620331
$arr[58,1] = -0.6323381214408973
$arr[59,0] = -0.2481888285055462
$arr[59,1] = -0.3683235323715158
$arr[60,0] = -0.1584826003877329
$arr[60,1] = -0.04380153611314741
$arr[61,0] = -0.3427909322360576
$arr[61,1] = -1.138481162617415
$arr[62,0] = -0.2617870985452501
$arr[62,1] = -0.6580405233340705
$arr[63,0] = -0.1909051150815745
$arr[63,1] = -0.6625356088187232
$arr[64,0] = -0.1242398638569124
$arr[64,1] = -0.1100233957378652
$arr[65,0] = -0.07863870715523616
$arr[65,1] = -0.2828814400549361

$ws.Range("A2:B67").Value = $arr
